$wb = $excel.ActiveWorkbook

# --- Worksheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")

# Row 19
$ws.Range("H19").Value = 677.4091
$ws.Range("I19").Value = 964.6667
$ws.Range("J19").Value = 478.53845
$ws.Range("K19").Value = 964.6667
$ws.Range("L19").Value = 478.53845
$ws.Range("M19").Value = -789.6667
$ws.Range("N19").Value = -828.53845

# Row 33
$ws.Range("H33").Value = 893.67645
$ws.Range("I33").Value = 729.5
$ws.Range("J33").Value = 1287.7
$ws.Range("K33").Value = 729.5
$ws.Range("L33").Value = 1287.7
$ws.Range("M33").Value = -500.5
$ws.Range("N33").Value = -1745.7

# Row 51
$ws.Range("H51").Value = 13890464
$ws.Range("I51").Value = 900
$ws.Range("J51").Value = 15874688
$ws.Range("K51").Value = 900
$ws.Range("L51").Value = 15874688
$ws.Range("M51").Value = -416
$ws.Range("N51").Value = -15875656

# Row 76
$ws.Range("H76").Value = 188876.33
$ws.Range("I76").Value = 446862.4
$ws.Range("J76").Value = 4600.5713
$ws.Range("K76").Value = 446862.4
$ws.Range("L76").Value = 4600.5713
$ws.Range("M76").Value = -446547.4
$ws.Range("N76").Value = -5230.5713

# Row 79
$ws.Range("H79").Value = 188876.33
$ws.Range("I79").Value = 446862.4
$ws.Range("J79").Value = 4600.5713
$ws.Range("K79").Value = 446862.4
$ws.Range("L79").Value = 4600.5713
$ws.Range("M79").Value = -445770.4
$ws.Range("N79").Value = -6784.5713

# Row 80
$ws.Range("H80").Value = 1366.24
$ws.Range("I80").Value = 1051.6765
$ws.Range("J80").Value = 2034.6875
$ws.Range("K80").Value = 3155.0295
$ws.Range("L80").Value = 6104.0625
$ws.Range("M80").Value = -2157.0295
$ws.Range("N80").Value = -8100.0625

# Row 83
$ws.Range("H83").Value = 1366.24
$ws.Range("I83").Value = 1051.6765
$ws.Range("J83").Value = 2034.6875
$ws.Range("K83").Value = 9465.0885
$ws.Range("L83").Value = 18312.1875
$ws.Range("M83").Value = -4473.0885
$ws.Range("N83").Value = -28296.1875

# Row 132
$ws.Range("H132").Value = 4807.5654
$ws.Range("I132").Value = 4889.727
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 14669.181
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -12139.181
$ws.Range("N132").Value = -14060

# Row 134
$ws.Range("H134").Value = 69933.336
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 69933.336
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 69933.336
$ws.Range("N134").Value = -80073.336

# Row 135
$ws.Range("H135").Value = 1057.5
$ws.Range("I135").Value = 1057.5
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 9517.5
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -6982.5

# Row 136
$ws.Range("H136").Value = 49000
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 49000
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 49000
$ws.Range("N136").Value = -59200

# Row 137
$ws.Range("H137").Value = 1446
$ws.Range("I137").Value = 1118.2084
$ws.Range("J137").Value = 2569.8572
$ws.Range("K137").Value = 3354.6252
$ws.Range("L137").Value = 7709.571599999999
$ws.Range("M137").Value = -804.6251999999999
$ws.Range("N137").Value = -12809.5716

# Row 139
$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

# Row 140
$ws.Range("H140").Value = 0
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()

# Row 141
$ws.Range("H141").Value = 8405
$ws.Range("I141").Value = 2680
$ws.Range("J141").Value = 21000
$ws.Range("K141").Value = 8040
$ws.Range("L141").Value = 63000
$ws.Range("M141").Value = -2860
$ws.Range("N141").Value = -73360

# --- Worksheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")

# Row 25
$ws.Range("H25").Value = 760.3333
$ws.Range("I25").Value = 760.3333
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 760.3333
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = -358.3333
$ws.Range("N25").ClearContents()

# Row 31
$ws.Range("H31").Value = 5504.2
$ws.Range("I31").Value = 5504.2
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 5504.2
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -5210.2

# Row 32
$ws.Range("H32").Value = 13256.83
$ws.Range("I32").Value = 3056.658
$ws.Range("J32").Value = 31714.285
$ws.Range("K32").Value = 3056.658
$ws.Range("L32").Value = 31714.285
$ws.Range("M32").Value = -2769.658
$ws.Range("N32").Value = -32288.285

# Row 35
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("M35").ClearContents()

# Row 41
$ws.Range("H41").Value = 5991
$ws.Range("I41").Value = 1056
$ws.Range("J41").Value = 7224.75
$ws.Range("K41").Value = 1056
$ws.Range("L41").Value = 7224.75
$ws.Range("M41").Value = -642
$ws.Range("N41").Value = -8052.75

# --- Worksheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")

# Row 12
$ws.Range("H12").Value = 724.875
$ws.Range("I12").Value = 299.83334
$ws.Range("J12").Value = 2000
$ws.Range("K12").Value = 299.83334
$ws.Range("L12").Value = 2000
$ws.Range("M12").Value = -131.83334
$ws.Range("N12").Value = -2336

# Row 24
$ws.Range("H24").Value = 369.6154
$ws.Range("I24").Value = 466.7
$ws.Range("J24").Value = 46
$ws.Range("K24").Value = 466.7
$ws.Range("L24").Value = 46
$ws.Range("M24").Value = -231.7
$ws.Range("N24").Value = -516

# Row 37
$ws.Range("H37").Value = 735
$ws.Range("I37").Value = 650.2857
$ws.Range("J37").Value = 932.6667
$ws.Range("K37").Value = 650.2857
$ws.Range("L37").Value = 932.6667
$ws.Range("M37").Value = -513.2857
$ws.Range("N37").Value = -1206.6667

# Row 82
$ws.Range("H82").Value = 18253.5
$ws.Range("I82").Value = 18253.5
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 18253.5
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -17870.5

# Row 85
$ws.Range("H85").Value = 18253.5
$ws.Range("I85").Value = 18253.5
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 18253.5
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -16927.5

# Row 86
$ws.Range("H86").Value = 16669159
$ws.Range("I86").Value = 20002190
$ws.Range("J86").Value = 4003.5
$ws.Range("K86").Value = 20002190
$ws.Range("L86").Value = 4003.5
$ws.Range("M86").Value = -20001067
$ws.Range("N86").Value = -6249.5

# Row 89
$ws.Range("H89").Value = 16669159
$ws.Range("I89").Value = 20002190
$ws.Range("J89").Value = 4003.5
$ws.Range("K89").Value = 100010950
$ws.Range("L89").Value = 20017.5
$ws.Range("M89").Value = -100005334
$ws.Range("N89").Value = -31249.5

# Row 97
$ws.Range("H97").Value = 17520
$ws.Range("I97").Value = 3300
$ws.Range("J97").Value = 27000
$ws.Range("K97").Value = 3300
$ws.Range("L97").Value = 27000
$ws.Range("M97").Value = -2309
$ws.Range("N97").Value = -28982

# --- Worksheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")

# Row 16
$ws.Range("H16").Value = 35716908
$ws.Range("I16").Value = 50002300
$ws.Range("J16").Value = 3425
$ws.Range("K16").Value = 50002300
$ws.Range("L16").Value = 3425
$ws.Range("M16").Value = -50002013
$ws.Range("N16").Value = -3999

# Row 113
$ws.Range("H113").Value = 35716908
$ws.Range("I113").Value = 50002300
$ws.Range("J113").Value = 3425
$ws.Range("K113").Value = 50002300
$ws.Range("L113").Value = 3425
$ws.Range("M113").Value = -50000130
$ws.Range("N113").Value = -7765

# --- Worksheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")

# Row 5
$ws.Range("H5").Value = 592
$ws.Range("I5").Value = 602.2222
$ws.Range("J5").Value = 500
$ws.Range("K5").Value = 1806.6666
$ws.Range("L5").Value = 1500
$ws.Range("M5").Value = -1694.6666
$ws.Range("N5").Value = -1724

# Row 9
$ws.Range("H9").Value = 39481.965
$ws.Range("I9").Value = 122822.445
$ws.Range("J9").Value = 3764.6191
$ws.Range("K9").Value = 368467.335
$ws.Range("L9").Value = 11293.8573
$ws.Range("M9").Value = -368243.335
$ws.Range("N9").Value = -11741.8573

# Row 62
$ws.Range("H62").Value = 5197.8887
$ws.Range("I62").Value = 3912
$ws.Range("J62").Value = 5358.625
$ws.Range("K62").Value = 11736
$ws.Range("L62").Value = 16075.875
$ws.Range("M62").Value = -11050
$ws.Range("N62").Value = -17447.875

# Row 65
$ws.Range("H65").Value = 5197.8887
$ws.Range("I65").Value = 3912
$ws.Range("J65").Value = 5358.625
$ws.Range("K65").Value = 35208
$ws.Range("L65").Value = 48227.625
$ws.Range("M65").Value = -31776
$ws.Range("N65").Value = -55091.625

# Row 113
$ws.Range("H113").Value = 722.8261
$ws.Range("I113").Value = 488.69565
$ws.Range("J113").Value = 956.95654
$ws.Range("K113").Value = 1466.08695
$ws.Range("L113").Value = 2870.86962
$ws.Range("M113").Value = 703.9130500000001
$ws.Range("N113").Value = -7210.869619999999

# Row 131
$ws.Range("H131").Value = 883.24
$ws.Range("I131").Value = 415.3846
$ws.Range("J131").Value = 953.1494
$ws.Range("K131").Value = 1246.1538
$ws.Range("L131").Value = 2859.4482
$ws.Range("M131").Value = 3793.8462
$ws.Range("N131").Value = -12939.4482

# Row 135
$ws.Range("H135").Value = 592
$ws.Range("I135").Value = 602.2222
$ws.Range("J135").Value = 500
$ws.Range("K135").Value = 5419.999800000001
$ws.Range("L135").Value = 4500
$ws.Range("M135").Value = -2884.999800000001
$ws.Range("N135").Value = -9570

# --- Worksheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")

# Row 41
$ws.Range("H41").Value = 6450.3335
$ws.Range("I41").Value = 351
$ws.Range("J41").Value = 9500
$ws.Range("K41").Value = 351
$ws.Range("L41").Value = 9500
$ws.Range("M41").Value = 4
$ws.Range("N41").Value = -10210

# Row 57
$ws.Range("H57").Value = 13061
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 13061
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 13061
$ws.Range("N57").Value = -14701
$ws.Range("M57").ClearContents()

# Row 113
$ws.Range("H113").Value = 2005.5454
$ws.Range("I113").Value = 1712.2
$ws.Range("J113").Value = 2250
$ws.Range("K113").Value = 1712.2
$ws.Range("L113").Value = 2250
$ws.Range("M113").Value = 457.8
$ws.Range("N113").Value = -6590
